$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4104) | Out-Null  # xlPasteAll
$excel.CutCopyMode = 0
